# Swap B, D, E, F, G values between paired rows (reversing a prior swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(149,150),
    @(313,314),
    @(316,318),
    @(346,347),
    @(350,351),
    @(382,383),
    @(389,390),
    @(400,401),
    @(457,458),
    @(536,537),
    @(586,587),
    @(590,591),
    @(593,594),
    @(601,602),
    @(709,710),
    @(715,716),
    @(720,721)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
